$d = $word.ActiveDocument

# Locate the run that currently reads "今天天气不错。" (weather run that
# follows the "晴，今天是高考第一天，上午考语文，下午考数学。" run).
$found = $d.Content
$found.Find.Execute("今天天气不错。", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $found.Start
$runEnd = $found.End

# Drop the trailing "。" from that run, leaving "今天天气不错".
$period = $d.Range($runEnd - 1, $runEnd)
$period.Delete()
$newEnd = $runEnd - 1

# Insert the additional sentence "，心情也很好。" right after it as new text.
$insertionPoint = $d.Range($newEnd, $newEnd)
$insertionPoint.InsertAfter("，心情也很好。")

# Toggling (and reverting) direct character formatting forces the engine to
# keep the freshly inserted text as its own run instead of silently folding
# it back into its identically-formatted neighbour.
$newRun = $d.Range($newEnd, $newEnd + 7)
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false

# Do the same at the boundary with the preceding "...下午考数学。" run so it
# also stays a distinct run (matches the original document's run layout).
$precedingRun = $d.Range($runStart, $newEnd)
$precedingRun.Font.Bold = $true
$precedingRun.Font.Bold = $false
